$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 5000120
$ws.Range("I5").Value = 5000120
$ws.Range("K5").Value = 5000120
$ws.Range("M5").Value = -5000005
$ws.Range("H64").Value = 3900.7222
$ws.Range("I64").Value = 3994
$ws.Range("J64").Value = 3889.0625
$ws.Range("K64").Value = 3994
$ws.Range("L64").Value = 3889.0625
$ws.Range("M64").Value = -3746
$ws.Range("N64").Value = -4385.0625
$ws.Range("H67").Value = 3900.7222
$ws.Range("I67").Value = 3994
$ws.Range("J67").Value = 3889.0625
$ws.Range("K67").Value = 3994
$ws.Range("L67").Value = 3889.0625
$ws.Range("M67").Value = -3136
$ws.Range("N67").Value = -5605.0625
$ws.Range("H106").Value = 5356.174
$ws.Range("I106").Value = 5356.174
$ws.Range("K106").Value = 5356.174
$ws.Range("M106").Value = -4725.174
$ws.Range("H139").Value = 106666
$ws.Range("J139").Value = 106666
$ws.Range("L139").Value = 106666
$ws.Range("N139").Value = -116946

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23047.18
$ws.Range("I32").Value = 13407.895
$ws.Range("K32").Value = 13407.895
$ws.Range("M32").Value = -13120.895
$ws.Range("H110").Value = 2946.6365
$ws.Range("I110").Value = 2946.6365
$ws.Range("K110").Value = 2946.6365
$ws.Range("M110").Value = -901.6365000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1607.2916
$ws.Range("I20").Value = 1566.0555
$ws.Range("J20").Value = 1731
$ws.Range("K20").Value = 1566.0555
$ws.Range("L20").Value = 1731
$ws.Range("M20").Value = -1319.0555
$ws.Range("N20").Value = -2225

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 170
$ws.Range("I4").Value = 129.23077
$ws.Range("K4").Value = 129.23077
$ws.Range("M4").Value = -17.23077000000001
$ws.Range("H31").Value = 3419
$ws.Range("I31").Value = 3377.75
$ws.Range("J31").Value = 3749
$ws.Range("K31").Value = 3377.75
$ws.Range("L31").Value = 3749
$ws.Range("M31").Value = -3082.75
$ws.Range("N31").Value = -4339
$ws.Range("H34").Value = 3419
$ws.Range("I34").Value = 3377.75
$ws.Range("J34").Value = 3749
$ws.Range("K34").Value = 3377.75
$ws.Range("L34").Value = 3749
$ws.Range("M34").Value = -3175.75
$ws.Range("N34").Value = -4153
$ws.Range("H35").Value = 1766.5834
$ws.Range("I35").Value = 1766.5834
$ws.Range("K35").Value = 1766.5834
$ws.Range("M35").Value = -1472.5834
$ws.Range("H54").Value = 39999.668
$ws.Range("J54").Value = 39999.668
$ws.Range("L54").Value = 39999.668
$ws.Range("N54").Value = -41315.668
$ws.Range("H69").Value = 39333
$ws.Range("I69").Value = 38999.5
$ws.Range("K69").Value = 38999.5
$ws.Range("M69").Value = -38250.5
$ws.Range("H72").Value = 39333
$ws.Range("I72").Value = 38999.5
$ws.Range("K72").Value = 116998.5
$ws.Range("M72").Value = -113254.5
$ws.Range("H122").Value = 3583.423
$ws.Range("I122").Value = 3590.8462
$ws.Range("K122").Value = 10772.5386
$ws.Range("M122").Value = -8322.5386
$ws.Range("H141").Value = 312247.66
$ws.Range("J141").Value = 312247.66
$ws.Range("L141").Value = 312247.66
$ws.Range("N141").Value = -322607.66

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 15804.909
$ws.Range("I139").Value = 13555
$ws.Range("J139").Value = 17090.572
$ws.Range("K139").Value = 40665
$ws.Range("L139").Value = 51271.716
$ws.Range("M139").Value = -35525
$ws.Range("N139").Value = -61551.716

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6404.1665
$ws.Range("J70").Value = 6922.4443
$ws.Range("L70").Value = 6922.4443
$ws.Range("N70").Value = -7462.4443
$ws.Range("H73").Value = 6404.1665
$ws.Range("J73").Value = 6922.4443
$ws.Range("L73").Value = 6922.4443
$ws.Range("N73").Value = -8794.444299999999
$ws.Range("H80").Value = 12499.75
$ws.Range("I80").Value = 6199.6
$ws.Range("K80").Value = 6199.6
$ws.Range("M80").Value = -5201.6
$ws.Range("H83").Value = 12499.75
$ws.Range("I83").Value = 6199.6
$ws.Range("K83").Value = 30998
$ws.Range("M83").Value = -26006
$ws.Range("H113").Value = 1519.6666
$ws.Range("I113").Value = 1547.6316
$ws.Range("K113").Value = 1547.6316
$ws.Range("M113").Value = 622.3684000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1115.4
$ws.Range("I16").Value = 1146.75
$ws.Range("K16").Value = 1146.75
$ws.Range("M16").Value = -976.75
$ws.Range("H22").Value = 734.8125
$ws.Range("I22").Value = 698.5714
$ws.Range("K22").Value = 698.5714
$ws.Range("M22").Value = -403.5714
$ws.Range("H27").Value = 734.8125
$ws.Range("I27").Value = 698.5714
$ws.Range("K27").Value = 698.5714
$ws.Range("M27").Value = -591.5714
$ws.Range("H32").Value = 2058.1667
$ws.Range("I32").Value = 2058.1667
$ws.Range("K32").Value = 2058.1667
$ws.Range("M32").Value = -1741.1667
$ws.Range("H82").Value = 2639.4
$ws.Range("I82").Value = 2248.75
$ws.Range("J82").Value = 4202
$ws.Range("K82").Value = 2248.75
$ws.Range("L82").Value = 4202
$ws.Range("M82").Value = -1887.75
$ws.Range("N82").Value = -4924
$ws.Range("H85").Value = 2639.4
$ws.Range("I85").Value = 2248.75
$ws.Range("J85").Value = 4202
$ws.Range("K85").Value = 2248.75
$ws.Range("L85").Value = 4202
$ws.Range("M85").Value = -1000.75
$ws.Range("N85").Value = -6698
$ws.Range("H122").Value = 4952.871
$ws.Range("I122").Value = 3507.5
$ws.Range("J122").Value = 5299.76
$ws.Range("K122").Value = 10522.5
$ws.Range("L122").Value = 15899.28
$ws.Range("M122").Value = -8072.5
$ws.Range("N122").Value = -20799.28

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 929.7778
$ws.Range("I100").Value = 955.17645
$ws.Range("K100").Value = 1910.3529
$ws.Range("M100").Value = -1369.3529
$ws.Range("H113").Value = 810.6
$ws.Range("I113").Value = 704.41174
$ws.Range("K113").Value = 2113.23522
$ws.Range("M113").Value = 56.76477999999997
$ws.Range("H122").Value = 4620.5
$ws.Range("I122").Value = 6274
$ws.Range("J122").Value = 2415.8333
$ws.Range("K122").Value = 18822
$ws.Range("L122").Value = 7247.499899999999
$ws.Range("M122").Value = -16372
$ws.Range("N122").Value = -12147.4999
$ws.Range("H135").Value = 66732.25
$ws.Range("J135").Value = 66732.25
$ws.Range("L135").Value = 66732.25
$ws.Range("N135").Value = -76872.25
